# Update database schema: remove the now-unused "time_stop" mandatory field
# from Sheet1 of the workbook. This deletes the entire row that holds the
# "time_stop" label (and its mirrored values in columns A-C), shifting all
# following rows up by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Find the row that contains "time_stop" in column A and delete it entirely
# so every cell below shifts up (rows, styles and all).
$timeStopRow = $ws.Range("A1:A20").Find("time_stop")
if ($timeStopRow -ne $null) {
    $ws.Rows.Item($timeStopRow.Row).Delete()
}

# Refresh the view to match the edited workbook: normal (non split) zoom,
# and move the active selection down to an empty row beneath the table.
$ws.Activate()
$excel.ActiveWindow.Zoom = 130
$ws.Range("B20").Select()
